# schedule.xlsx update — "update schedule and 解答"
#
# Semantic changes:
#   1. Row 7 (the "Java10/RecyclerView" task row) gets its completion
#      date filled in: C7 = 2018-04-18 (serial 43208), formatted like
#      the other date cells in column B/C.
#   2. A brand-new task row (row 8) is added:
#        A8 = new task description
#        B8 = expected date 2018-04-20 (serial 43210)
#        D8 = remark/notes text
#      (C8 / E8 are left blank, same as row 7 was before it was marked
#      complete.)
#   3. The active selection moves to D8 (where the user was last typing).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- 1. Fill in the completion date for row 7 -----------------------------
$ws.Cells.Item(7, 3).Value = 43208
$ws.Cells.Item(7, 3).NumberFormat = "m/d/yy"

# --- 2. Add the new row of data -------------------------------------------

# A8: task text — copy formatting (wrap text, border, font) from the cell
# directly above it so it matches the rest of the "task" column.
$ws.Cells.Item(7, 1).Copy() | Out-Null
$ws.Cells.Item(8, 1).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(8, 1).Value = "1.Java另一本书的第10章集合类" + [char]10 + "2.使用RecyclerView写一个仿微信列表(界面如:列表.jpg)"

# B8: expected date — copy formatting from the cell above (date number format).
$ws.Cells.Item(7, 2).Copy() | Out-Null
$ws.Cells.Item(8, 2).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(8, 2).Value = 43210

# D8: remark text — copy formatting from the cell above (wrap text, border, font).
$ws.Cells.Item(7, 4).Copy() | Out-Null
$ws.Cells.Item(8, 4).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(8, 4).Value = "看另一本Java从入门到精通。" + [char]10 + "界面只要关注红色框框部分，就是只关注列表的实现，另外有个黑框圈起来的可以忽略那个部分，也就是每个item大概由四个部分组成(给两天时间)"

$excel.CutCopyMode = 0

# --- 3. Move the selection to where the user left off ----------------------
$ws.Range("D8").Select() | Out-Null
